$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct a tiny floating point drift in the existing A9 timestamp
$ws.Range("A9").Value = 45809.39161978009

# Append the new price record as row 10
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat()
$ws.Range("A10").Value = 45810.39392270573
$ws.Range("B10").Value = "EVOWHEY PROTEIN"
$ws.Range("C10").Value = "2Kg"
$ws.Range("D10").Value = "34,90€"
